$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01762077648523471
$ws.Range("C2").Value = 0.2606746064175882
$ws.Range("D2").Value = 0.2060894145119571
$ws.Range("E2").Value = 0.4539707198839559
$ws.Range("F2").Value = 0.4580543675004212

$ws.Range("B3").Value = -0.03991108385890063
$ws.Range("C3").Value = 0.4106880171616463
$ws.Range("D3").Value = 0.3811762225484371
$ws.Range("E3").Value = 0.6173947056368698
$ws.Range("F3").Value = 0.6222338744333338

$ws.Range("B4").Value = 0.07916732041434071
$ws.Range("C4").Value = 0.3709301532769832
$ws.Range("D4").Value = 0.3042748425083037
$ws.Range("E4").Value = 0.5516111334158365
$ws.Range("F4").Value = 0.5514427994895704

$ws.Range("B5").Value = 0.01404048331219524
$ws.Range("C5").Value = 0.3932946038540242
$ws.Range("D5").Value = 0.3438932018012958
$ws.Range("E5").Value = 0.5864240801683503
$ws.Range("F5").Value = 0.5923313273423136
$ws.Range("G5").Value = 49

$ws.Range("B6").Value = 0.1396864754206484
$ws.Range("C6").Value = 0.3647466888290081
$ws.Range("D6").Value = 0.3239417985499088
$ws.Range("E6").Value = 0.5691588517715497
$ws.Range("F6").Value = 0.557590091612993
$ws.Range("G6").Value = 48

$ws.Range("B7").Value = 0.06803994534832956
$ws.Range("C7").Value = 0.3466446476016687
$ws.Range("D7").Value = 0.2880673810335034
$ws.Range("E7").Value = 0.5367190894998085
$ws.Range("F7").Value = 0.539348533148276
$ws.Range("G7").Value = 39

$ws.Range("B8").Value = 0.1037691244471713
$ws.Range("C8").Value = 0.330411853491841
$ws.Range("D8").Value = 0.2727802349308482
$ws.Range("E8").Value = 0.5222836728549421
$ws.Range("F8").Value = 0.5187423393692384
$ws.Range("G8").Value = 38

$ws.Range("B9").Value = 0.125076358035327
$ws.Range("C9").Value = 0.3482004317699195
$ws.Range("D9").Value = 0.3553016602029547
$ws.Range("E9").Value = 0.5960718582544848
$ws.Range("F9").Value = 0.5971938069896171
$ws.Range("G9").Value = 21

$ws.Range("B10").Value = 0.2709217060049215
$ws.Range("C10").Value = 0.3713999088374357
$ws.Range("D10").Value = 0.45908714801176
$ws.Range("E10").Value = 0.6775597006993258
$ws.Range("F10").Value = 0.6444819076758775
$ws.Range("G10").Value = 14

$ws.Range("B11").Value = 0.2443462408704248
$ws.Range("C11").Value = 0.2611723105820856
$ws.Range("D11").Value = 0.1488287689115677
$ws.Range("E11").Value = 0.3857833186019942
$ws.Range("F11").Value = 0.3337732828658925
